$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.99999994583665
$ws.Range("E2").Value = 0.99999994583665

# Row 3
$ws.Range("C3").Value = $true
$ws.Range("D3").Value = 0.05172475832356097
$ws.Range("E3").Value = 0.05172475832356097

# Row 4
$ws.Range("D4").Value = 0.01427608004938534
$ws.Range("E4").Value = 0.01427608004938534

# Row 5
$ws.Range("D5").Value = [double]"6.568755928308975E-26"
$ws.Range("E5").Value = [double]"6.568755928308975E-26"

# Row 6
$ws.Range("D6").Value = [double]"1.088997538693839E-15"
$ws.Range("E6").Value = [double]"1.088997538693839E-15"

# Row 7
$ws.Range("D7").Value = 0.9999578544772678
$ws.Range("E7").Value = [double]"4.21455227321843E-05"

# Row 8
$ws.Range("D8").Value = 0.9999999999870066
$ws.Range("E8").Value = [double]"1.299338414639806E-11"

# Row 10
$ws.Range("D10").Value = 0.914316228675556
$ws.Range("E10").Value = 0.08568377132444405

# Row 11
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1.68883740901947
$ws.Range("G11").Value = 0.9
